$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2-97) forward by 7 days,
# keeping the same time-of-day fraction (new week of data).
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 7
}

# Updated Actual Production (MW) values for the new week in column B
# (rows 25-38 get new non-zero readings, rows 39-66 reset to 0).
$bUpdates = @{
    25 = 4
    26 = 17
    27 = 40
    28 = 74
    29 = 127
    30 = 189
    31 = 259
    32 = 317
    33 = 381
    34 = 454
    35 = 531
    36 = 574
    37 = 659
    38 = 670
    39 = 0
    40 = 0
    41 = 0
    42 = 0
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 0
    49 = 0
    50 = 0
    51 = 0
    52 = 0
    53 = 0
    54 = 0
    55 = 0
    56 = 0
    57 = 0
    58 = 0
    59 = 0
    60 = 0
    61 = 0
    62 = 0
    63 = 0
    64 = 0
    65 = 0
    66 = 0
}

foreach ($r in $bUpdates.Keys) {
    $ws.Cells.Item($r, 2).Value = $bUpdates[$r]
}

Write-Host "Edit applied"
